$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-15 (columns A-F)
$data = @(
    @(45918, 5101.32827217152, 4933.43939847377, 7056, 6512.511375, -29.6407291124064),
    @(45919, 5103.30962648907, 4197.493965356, 2952, 6516.650959, 110.784804077789),
    @(45920, 1207.12767057941, 2143.11766689043, 2952, 2319.790459, 12.6575189712921),
    @(45921, 1096.63641792581, 2118.00850779633, 2952, 2202.28905, 11.3192141612716),
    @(45922, 5049.38273576318, 4726.8253508235, 2952, 6710.218759, 143.152557252514),
    @(45923, 5230.43751372964, 5024.02246794799, 2952, 6710.218759, 147.991821384098),
    @(45924, 5230.43751372964, 5049.57230201316, 2952, 6710.218759, 149.05639780348),
    @(45925, 5230.43751372964, 5726.1675249974, 2952, 6710.218759, 177.247865427824),
    @(45926, 5230.43751372964, 5579.77902050718, 2952, 6710.218759, 171.148344407397),
    @(45927, 1198.04912219854, 3643.47003663811, 2952, 2321.483151, 75.6210027266487),
    @(45928, 1070.29084929054, 3604.19170099456, 2952, 2183.373747, 73.5531082793344),
    @(45929, 5457.44872542722, 6078.26410629411, 2952, 7031.626518, 195.85174578612),
    @(45930, 5457.44872542722, 6090.41862578345, 2952, 7031.626518, 196.358184098176),
    @(45931, 4260.05383201679, 5810.81868970956, 3692, 6097.99183, 164.864861987199)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 0; $col -lt $values.Length; $col++) {
        $ws.Cells.Item($row, $col + 1).Value = $values[$col]
    }
}
